$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 161114.84   # H15
$ws.Cells.Item(15, 9).Value = 161114.84   # I15
$ws.Cells.Item(15, 11).Value = 483344.52   # K15
$ws.Cells.Item(15, 13).Value = -483175.52   # M15
$ws.Cells.Item(43, 8).Value = 8668.655000000001   # H43
$ws.Cells.Item(43, 9).Value = 6584   # I43
$ws.Cells.Item(43, 10).Value = 9212.478999999999   # J43
$ws.Cells.Item(43, 11).Value = 6584   # K43
$ws.Cells.Item(43, 12).Value = 9212.478999999999   # L43
$ws.Cells.Item(43, 13).Value = -6515   # M43
$ws.Cells.Item(43, 14).Value = -9350.478999999999   # N43
$ws.Cells.Item(53, 8).Value = 84193.75   # H53
$ws.Cells.Item(53, 9).Value = 200670.4   # I53
$ws.Cells.Item(53, 11).Value = 200670.4   # K53
$ws.Cells.Item(53, 13).Value = -200033.4   # M53
$ws.Cells.Item(113, 8).Value = 3847.8   # H113
$ws.Cells.Item(113, 9).Value = 5494.5   # I113
$ws.Cells.Item(113, 11).Value = 5494.5   # K113
$ws.Cells.Item(113, 13).Value = -2240.5   # M113
$ws.Cells.Item(137, 8).Value = 7667.263   # H137
$ws.Cells.Item(137, 9).Value = 8219.333000000001   # I137
$ws.Cells.Item(137, 10).Value = 6720.857   # J137
$ws.Cells.Item(137, 11).Value = 24657.999   # K137
$ws.Cells.Item(137, 12).Value = 20162.571   # L137
$ws.Cells.Item(137, 13).Value = -22107.999   # M137
$ws.Cells.Item(137, 14).Value = -25262.571   # N137
$ws.Cells.Item(138, 8).Value = 3654.2888   # H138
$ws.Cells.Item(138, 9).Value = 1511.2858   # I138
$ws.Cells.Item(138, 10).Value = 5529.4165   # J138
$ws.Cells.Item(138, 11).Value = 4533.857400000001   # K138
$ws.Cells.Item(138, 12).Value = 16588.2495   # L138
$ws.Cells.Item(138, 13).Value = 606.1425999999992   # M138
$ws.Cells.Item(138, 14).Value = -26868.2495   # N138

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(35, 8).Value = 14450   # H35
$ws.Cells.Item(35, 9).Value = 20000   # I35
$ws.Cells.Item(35, 10).Value = 8900   # J35
$ws.Cells.Item(35, 11).Value = 20000   # K35
$ws.Cells.Item(35, 12).Value = 8900   # L35
$ws.Cells.Item(35, 13).Value = -19594   # M35
$ws.Cells.Item(35, 14).Value = -9712   # N35
$ws.Cells.Item(39, 8).Value = 13609.333   # H39
$ws.Cells.Item(39, 9).Value = 1000   # I39
$ws.Cells.Item(39, 11).Value = 1000   # K39
$ws.Cells.Item(39, 13).Value = -480   # M39
$ws.Cells.Item(61, 8).Value = 4285.6665   # H61
$ws.Cells.Item(61, 9).Value = 3744.6428   # I61
$ws.Cells.Item(61, 10).Value = 5662.8184   # J61
$ws.Cells.Item(61, 11).Value = 3744.6428   # K61
$ws.Cells.Item(61, 12).Value = 5662.8184   # L61
$ws.Cells.Item(61, 13).Value = -3532.6428   # M61
$ws.Cells.Item(61, 14).Value = -6086.8184   # N61
$ws.Cells.Item(74, 8).Value = 5084.45   # H74
$ws.Cells.Item(74, 9).Value = 6144.4546   # I74
$ws.Cells.Item(74, 10).Value = 3788.889   # J74
$ws.Cells.Item(74, 11).Value = 6144.4546   # K74
$ws.Cells.Item(74, 12).Value = 3788.889   # L74
$ws.Cells.Item(74, 13).Value = -5270.4546   # M74
$ws.Cells.Item(74, 14).Value = -5536.889   # N74
$ws.Cells.Item(76, 8).Value = 32349.5   # H76
$ws.Cells.Item(76, 10).Value = 32349.5   # J76
$ws.Cells.Item(76, 12).Value = 32349.5   # L76
$ws.Cells.Item(76, 14).Value = -33025.5   # N76
$ws.Cells.Item(77, 8).Value = 5084.45   # H77
$ws.Cells.Item(77, 9).Value = 6144.4546   # I77
$ws.Cells.Item(77, 10).Value = 3788.889   # J77
$ws.Cells.Item(77, 11).Value = 30722.273   # K77
$ws.Cells.Item(77, 12).Value = 18944.445   # L77
$ws.Cells.Item(77, 13).Value = -26354.273   # M77
$ws.Cells.Item(77, 14).Value = -27680.445   # N77
$ws.Cells.Item(79, 8).Value = 32349.5   # H79
$ws.Cells.Item(79, 10).Value = 32349.5   # J79
$ws.Cells.Item(79, 12).Value = 32349.5   # L79
$ws.Cells.Item(79, 14).Value = -34689.5   # N79
$ws.Cells.Item(136, 8).Value = 4285.6665   # H136
$ws.Cells.Item(136, 9).Value = 3744.6428   # I136
$ws.Cells.Item(136, 10).Value = 5662.8184   # J136
$ws.Cells.Item(136, 11).Value = 11233.9284   # K136
$ws.Cells.Item(136, 12).Value = 16988.4552   # L136
$ws.Cells.Item(136, 13).Value = -8683.928400000001   # M136
$ws.Cells.Item(136, 14).Value = -22088.4552   # N136

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(8, 8).Value = 4501.6665   # H8
$ws.Cells.Item(8, 9).Value = 1000   # I8
$ws.Cells.Item(8, 10).Value = 6252.5   # J8
$ws.Cells.Item(8, 11).Value = 1000   # K8
$ws.Cells.Item(8, 12).Value = 6252.5   # L8
$ws.Cells.Item(8, 13).Value = -860   # M8
$ws.Cells.Item(8, 14).Value = -6532.5   # N8
$ws.Cells.Item(134, 8).Value = 5573932   # H134
$ws.Cells.Item(134, 9).Value = 6267236   # I134
$ws.Cells.Item(134, 11).Value = 18801708   # K134
$ws.Cells.Item(134, 13).Value = -18799173   # M134

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5886.476   # H31
$ws.Cells.Item(31, 9).Value = 1721   # I31
$ws.Cells.Item(31, 11).Value = 1721   # K31
$ws.Cells.Item(31, 13).Value = -1426   # M31
$ws.Cells.Item(34, 8).Value = 5886.476   # H34
$ws.Cells.Item(34, 9).Value = 1721   # I34
$ws.Cells.Item(34, 11).Value = 1721   # K34
$ws.Cells.Item(34, 13).Value = -1519   # M34
$ws.Cells.Item(99, 8).Value = 22225362   # H99
$ws.Cells.Item(99, 9).Value = 27780652   # I99
$ws.Cells.Item(99, 11).Value = 27780652   # K99
$ws.Cells.Item(99, 13).Value = -27779154   # M99
$ws.Cells.Item(126, 8).Value = 22225362   # H126
$ws.Cells.Item(126, 9).Value = 27780652   # I126
$ws.Cells.Item(126, 11).Value = 83341956   # K126
$ws.Cells.Item(126, 13).Value = -83339486   # M126
$ws.Cells.Item(132, 8).Value = 23694   # H132
$ws.Cells.Item(132, 9).Value = 22996.334   # I132
$ws.Cells.Item(132, 11).Value = 68989.00199999999   # K132
$ws.Cells.Item(132, 13).Value = -66459.00199999999   # M132

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1196.3334   # H68
$ws.Cells.Item(68, 10).Value = 1599   # J68
$ws.Cells.Item(68, 12).Value = 4797   # L68
$ws.Cells.Item(68, 14).Value = -6419   # N68
$ws.Cells.Item(71, 8).Value = 1196.3334   # H71
$ws.Cells.Item(71, 10).Value = 1599   # J71
$ws.Cells.Item(71, 12).Value = 14391   # L71
$ws.Cells.Item(71, 14).Value = -22503   # N71
$ws.Cells.Item(107, 8).Value = 4472.9697   # H107
$ws.Cells.Item(107, 10).Value = 4597.125   # J107
$ws.Cells.Item(107, 12).Value = 13791.375   # L107
$ws.Cells.Item(107, 14).Value = -17631.375   # N107
$ws.Cells.Item(120, 8).Value = 31474.875   # H120
$ws.Cells.Item(120, 10).Value = 49760   # J120
$ws.Cells.Item(120, 12).Value = 149280   # L120
$ws.Cells.Item(120, 14).Value = -158956   # N120
$ws.Cells.Item(132, 8).Value = 2978.1   # H132
$ws.Cells.Item(132, 10).Value = 4480.3335   # J132
$ws.Cells.Item(132, 12).Value = 40323.0015   # L132
$ws.Cells.Item(132, 14).Value = -45383.0015   # N132
$ws.Cells.Item(140, 8).Value = 41668400   # H140
$ws.Cells.Item(140, 9).Value = 47102890   # I140
$ws.Cells.Item(140, 11).Value = 141308670   # K140
$ws.Cells.Item(140, 13).Value = -141303490   # M140

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(13, 8).Value = 3491   # H13
$ws.Cells.Item(13, 9).Value = 4113.75   # I13
$ws.Cells.Item(13, 11).Value = 4113.75   # K13
$ws.Cells.Item(13, 13).Value = -3974.75   # M13
$ws.Cells.Item(42, 8).Value = 30000   # H42
$ws.Cells.Item(42, 10).Value = 30000   # J42
$ws.Cells.Item(42, 12).Value = 30000   # L42
$ws.Cells.Item(42, 14).Value = -30970   # N42
$ws.Cells.Item(115, 8).Value = 30000   # H115
$ws.Cells.Item(115, 10).Value = 30000   # J115
$ws.Cells.Item(115, 12).Value = 30000   # L115
$ws.Cells.Item(115, 14).Value = -32350   # N115
$ws.Cells.Item(132, 8).Value = 26320200   # H132
$ws.Cells.Item(132, 9).Value = 38466028   # I132
$ws.Cells.Item(132, 11).Value = 115398084   # K132
$ws.Cells.Item(132, 13).Value = -115395554   # M132

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5941.5713   # H7
$ws.Cells.Item(7, 9).Value = 3766.6365   # I7
$ws.Cells.Item(7, 11).Value = 3766.6365   # K7
$ws.Cells.Item(7, 13).Value = -3654.6365   # M7
$ws.Cells.Item(11, 8).Value = 11761297   # H11
$ws.Cells.Item(11, 9).Value = 12828869   # I11
$ws.Cells.Item(11, 10).Value = 18007   # J11
$ws.Cells.Item(11, 11).Value = 12828869   # K11
$ws.Cells.Item(11, 12).Value = 18007   # L11
$ws.Cells.Item(11, 13).Value = -12828729   # M11
$ws.Cells.Item(11, 14).Value = -18287   # N11
$ws.Cells.Item(13, 8).Value = 16679500   # H13
$ws.Cells.Item(13, 9).Value = 16679500   # I13
$ws.Cells.Item(13, 11).Value = 16679500   # K13
$ws.Cells.Item(13, 13).Value = -16679360   # M13
$ws.Cells.Item(17, 8).Value = 353491460   # H17
$ws.Cells.Item(17, 9).Value = 25015000   # I17
$ws.Cells.Item(17, 11).Value = 25015000   # K17
$ws.Cells.Item(17, 13).Value = -25014830   # M17
$ws.Cells.Item(46, 8).Value = 71429730   # H46
$ws.Cells.Item(46, 9).Value = 1350   # I46
$ws.Cells.Item(46, 10).Value = 166667570   # J46
$ws.Cells.Item(46, 11).Value = 1350   # K46
$ws.Cells.Item(46, 12).Value = 166667570   # L46
$ws.Cells.Item(46, 13).Value = -1162   # M46
$ws.Cells.Item(46, 14).Value = -166667946   # N46
$ws.Cells.Item(61, 8).Value = 5940.115   # H61
$ws.Cells.Item(61, 9).Value = 4656.1665   # I61
$ws.Cells.Item(61, 10).Value = 8829   # J61
$ws.Cells.Item(61, 11).Value = 4656.1665   # K61
$ws.Cells.Item(61, 12).Value = 8829   # L61
$ws.Cells.Item(61, 13).Value = -4454.1665   # M61
$ws.Cells.Item(61, 14).Value = -9233   # N61
$ws.Cells.Item(92, 8).Value = 15879.6   # H92
$ws.Cells.Item(92, 9).Value = 18000   # I92
$ws.Cells.Item(92, 10).Value = 15349.5   # J92
$ws.Cells.Item(92, 11).Value = 18000   # K92
$ws.Cells.Item(92, 12).Value = 15349.5   # L92
$ws.Cells.Item(92, 13).Value = -15504   # M92
$ws.Cells.Item(92, 14).Value = -20341.5   # N92
$ws.Cells.Item(113, 8).Value = 5940.115   # H113
$ws.Cells.Item(113, 9).Value = 4656.1665   # I113
$ws.Cells.Item(113, 10).Value = 8829   # J113
$ws.Cells.Item(113, 11).Value = 4656.1665   # K113
$ws.Cells.Item(113, 12).Value = 8829   # L113
$ws.Cells.Item(113, 13).Value = -2486.1665   # M113
$ws.Cells.Item(113, 14).Value = -13169   # N113
$ws.Cells.Item(126, 8).Value = 5941.5713   # H126
$ws.Cells.Item(126, 9).Value = 3766.6365   # I126
$ws.Cells.Item(126, 11).Value = 11299.9095   # K126
$ws.Cells.Item(126, 13).Value = -8829.9095   # M126
$ws.Cells.Item(132, 8).Value = 7507   # H132
$ws.Cells.Item(132, 9).Value = 8316.666999999999   # I132
$ws.Cells.Item(132, 10).Value = 6899.75   # J132
$ws.Cells.Item(132, 11).Value = 24950.001   # K132
$ws.Cells.Item(132, 12).Value = 20699.25   # L132
$ws.Cells.Item(132, 13).Value = -22420.001   # M132
$ws.Cells.Item(132, 14).Value = -25759.25   # N132

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 2012   # H14
$ws.Cells.Item(14, 9).Value = 2012   # I14
$ws.Cells.Item(14, 10).Value = 0   # J14
$ws.Cells.Item(14, 11).Value = 2012   # K14
$ws.Cells.Item(14, 12).Value = 0   # L14
$ws.Cells.Item(14, 13).Value = -1844   # M14
$ws.Cells.Item(14, 14).Value = $null   # N14 (delete)
$ws.Cells.Item(17, 8).Value = 25015000   # H17
$ws.Cells.Item(17, 9).Value = 25015000   # I17
$ws.Cells.Item(17, 11).Value = 25015000   # K17
$ws.Cells.Item(17, 13).Value = -25014828   # M17
$ws.Cells.Item(47, 8).Value = 30000   # H47
$ws.Cells.Item(47, 10).Value = 30000   # J47
$ws.Cells.Item(47, 12).Value = 30000   # L47
$ws.Cells.Item(47, 14).Value = -31144   # N47
$ws.Cells.Item(81, 8).Value = 1317.7333   # H81
$ws.Cells.Item(81, 10).Value = 2000   # J81
$ws.Cells.Item(81, 12).Value = 4000   # L81
$ws.Cells.Item(81, 14).Value = -6122   # N81
$ws.Cells.Item(84, 8).Value = 1317.7333   # H84
$ws.Cells.Item(84, 10).Value = 2000   # J84
$ws.Cells.Item(84, 12).Value = 20000   # L84
$ws.Cells.Item(84, 14).Value = -30608   # N84
$ws.Cells.Item(113, 8).Value = 7938365.5   # H113
$ws.Cells.Item(113, 9).Value = 12822365   # I113
$ws.Cells.Item(113, 11).Value = 38467095   # K113
$ws.Cells.Item(113, 13).Value = -38464925   # M113
